$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: D2 1000 -> 4000
$ws.Range("D2").Value = 4000

# Update row 3: A3 45889 -> 45890, C3 395 -> 394, D3 1000 -> 8000
$ws.Range("A3").Value = 45890
$ws.Range("C3").Value = 394
$ws.Range("D3").Value = 8000

# Clear row 4 data (A4, C4, D4) - keep A4's style
$ws.Range("A4").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("D4").ClearContents()

# Clear row 5 data (A5, C5, D5)
$ws.Range("A5").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("D5").ClearContents()

# Update selection to A3
$ws.Range("A3").Select()
